$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $values[$i]
    }
}

# --- Update the "last updated" timestamp (A1) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 14 de Mayo de 2020 a las 19:05"

# --- Plain value refresh (no rank change) ---
# Row 4: Estados Unidos
Set-Row 4 @(1438426, 8078, 311693, 1041034, 16337, 502, 85699)

# Row 11: Alemania
Set-Row 11 @(174478, 380, 150300, 16294, 1465, 23, 7884)

# Row 12: Turquia
Set-Row 12 @(144749, 1635, 104030, 36712, 963, 55, 4007)

# Row 15: India
Set-Row 15 @(81990, 3935, 27956, 51385, 0, 98, 2649)

# --- Egipto overtakes Dinamarca (rows 47/48 swap rank, row 47 gets fresh data) ---
$ws.Cells.Item(47, 1).Value = "Egipto"
Set-Row 47 @(10829, 398, 2626, 7632, 41, 15, 571)
$ws.Cells.Item(48, 1).Value = "Dinamarca"
Set-Row 48 @(10713, 46, 8805, 1371, 35, 4, 537)

# --- Camerun overtakes Azerbaiyan (rows 72/73 swap rank, row 72 gets fresh data) ---
$ws.Cells.Item(72, 1).Value = "Camerun"
Set-Row 72 @(2954, 154, 1553, 1262, 28, 3, 139)
$ws.Cells.Item(73, 1).Value = "Azerbaiyan"
Set-Row 73 @(2879, 121, 1833, 1011, 29, 0, 35)

# --- Benin jumps ahead of Congo / Isla de Man / Mauricio (rows 133-136 shift down, row 133 gets fresh data) ---
$ws.Cells.Item(133, 1).Value = "Benin"
Set-Row 133 @(339, 12, 83, 254, 0, 0, 2)
$ws.Cells.Item(134, 1).Value = "Congo"
Set-Row 134 @(333, 0, 53, 269, 0, 0, 11)
$ws.Cells.Item(135, 1).Value = "Isla de Man"
Set-Row 135 @(332, 0, 274, 35, 20, 0, 23)
$ws.Cells.Item(136, 1).Value = "Mauricio"
Set-Row 136 @(332, 0, 322, 0, 0, 0, 10)

# --- Belice overtakes Santa Lucia (rows 194/195 trade data exactly, no new values) ---
$ws.Cells.Item(194, 1).Value = "Belice"
Set-Row 194 @(18, 0, 16, 0, 0, 0, 2)
$ws.Cells.Item(195, 1).Value = "Santa Lucia"
Set-Row 195 @(18, 0, 18, 0, 0, 0, 0)

# --- San Bartolome overtakes Sahara Occidental (rows 215/216, identical data so no visible number change) ---
$ws.Cells.Item(215, 1).Value = "San Bartolome"
$ws.Cells.Item(216, 1).Value = "Sahara Occidental"
